# Generate Report for Handback
# Update the "Correspond Handoff Datetime" (col E) and "Correspond Handback DateTime"
# (col H) values for the e8579b9a... file row (row 3) on both the zh-cn and de-de
# handback-status worksheets, reflecting a newer handback report run.

$wb = $excel.ActiveWorkbook

# zh-cn sheet: row 3 is the e8579b9a...zh-cn.xlf file
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E3").Value = "2016-03-14 02:31:18"
$wsZhCn.Range("H3").Value = "2016-03-14 02:31:39"

# de-de sheet: row 3 is the e8579b9a...de-de.xlf file
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E3").Value = "2016-03-14 02:31:21"
$wsDeDe.Range("H3").Value = "2016-03-14 02:31:44"
